# Add new "multi-selection" localization entries to the FALanguage tag sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 599-604 were typed normally: column A (key), then column B (caption),
# row by row.
$normalRows = @(
    @("MultiSelectionOptionCaption",        "多选选项"),
    @("MultiSelectionOptionConnected",       "相连"),
    @("MultiSelectionOptionSameTileSet",     "相同地形组"),
    @("MultiSelectionOptionConsiderLAT",     "考虑LAT"),
    @("MultiSelectionOptionSameHeight",      "相同高度"),
    @("MultiSelectionOptionSameBaiscHeight", "相同基础高度")
)

$startRow = 599
for ($i = 0; $i -lt $normalRows.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $normalRows[$i][0]
    $ws.Cells.Item($row, 2).Value = $normalRows[$i][1]
}

# Rows 605-606 were filled out of the usual A-then-B order: both captions
# (column B) were entered first, then the keys (column A) were filled back
# in reverse row order. Reproduce that exact sequence so the shared-string
# table ends up in the same order as the authored workbook.
$ws.Cells.Item(605, 2).Value = "自定义添加"
$ws.Cells.Item(606, 2).Value = "自定义删除"
$ws.Cells.Item(606, 1).Value = "MultiSelectionCustomDelete"
$ws.Cells.Item(605, 1).Value = "MultiSelectionCustomAdd"

# Match the saved view state recorded in the workbook after the edit.
$ws.Application.ActiveWindow.ScrollRow = 579
$ws.Range("A601").Select()
